$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "[Marie%Lin%marilin@ms2.mmh.org.tw%1,                     Hsiang-Kuang%Tseng%drtseng@anet.net.tw%1,                     Jean A%Trejaut%jtrejaut@ms1.mmh.org.tw%1,                     Hui-Lin%Lee%leehl@ms1.mmh.org.tw%1,                     Jun-Hun%Loo%junhun@ms1.mmh.org.tw%1,                     Chen-Chung%Chu%chucc@ms1.mmh.org.tw%1,                     Pei-Jan%Chen%g8640@ms13.hinet.net%1,                     Ying-Wen%Su%yingwen_su@yahoo.com.tw%1,                     Ken Hong%Lim%khlim@seed.net.tw%1,                     Zen-Uong%Tsai%tsaizu@ms1.mmh.org.tw%1,                     Ruey-Yi%Lin%lin-8424@mail.hoping.gov.tw%1,                     Ruey-Shiung%Lin%rueyshiung@yahoo.com%1,                     Chun-Hsiung%Huang%chhuang@ms2.mmh.org.tw%1]"
$ws.Cells.Item(2, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(3, 5).Value = "[Rossa W K%Chiu%NULL%1,                     Nelson L S%Tang%NULL%1,                     David S C%Hui%NULL%1,                     Grace T Y%Chung%NULL%1,                     Stephen S C%Chim%NULL%1,                     K C Allen%Chan%NULL%1,                     Ying-man%Sung%NULL%1,                     Louis Y S%Chan%NULL%1,                     Yu-kwan%Tong%NULL%1,                     Wing-shan%Lee%NULL%1,                     Paul K S%Chan%NULL%1,                     Y M Dennis%Lo%loym@cuhk.edu.hk%1]"
$ws.Cells.Item(3, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(4, 5).Value = "[Satoru%Itoyama%NULL%2,                     Naoto%Keicho%NULL%3,                     Tran%Quy%NULL%4,                     Nguyen Chi%Phi%NULL%4,                     Hoang Thuy%Long%NULL%4,                     Le Dang%Ha%NULL%4,                     Vo Van%Ban%NULL%4,                     Jun%Ohashi%NULL%1,                     Minako%Hijikata%NULL%4,                     Ikumi%Matsushita%NULL%3,                     Akihiko%Kawana%NULL%1,                     Hideki%Yanai%NULL%3,                     Teruo%Kirikae%NULL%3,                     Tadatoshi%Kuratsuji%NULL%3,                     Takehiko%Sasazuki%NULL%4]"
$ws.Cells.Item(4, 9).Value = "_PMC_elsevier_CROSSREF"

$ws.Cells.Item(5, 5).Value = "[Margaret H. L.%Ng%margaretng@cuhk.edu.hk%1,                     Kin-Mang%Lau%NULL%1,                     Libby%Li%NULL%1,                     Suk-Hang%Cheng%NULL%1,                     Wing Y.%Chan%NULL%1,                     Pak K.%Hui%NULL%1,                     Benny%Zee%NULL%1,                     Chi-Bon%Leung%NULL%1,                     Joseph J. Y.%Sung%NULL%1]"
$ws.Cells.Item(5, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(6, 5).Value = "[KC Allen%Chan%allen@cuhk.edu.hk%1,                     Nelson LS%Tang%nelsontang@cuhk.edu.hk%1,                     David SC%Hui%b200945@mailserv.cuhk.edu.hk%1,                     Grace TY%Chung%gracechung@cuhk.edu.hk%1,                     Alan KL%Wu%alanklwu@yahoo.com%1,                     Stephen SC%Chim%sschim@cuhk.edu.hk%1,                     Rossa WK%Chiu%rossachiu@cuhk.edu.hk%1,                     Nelson%Lee%leelsn@yahoo.com%1,                     KW%Choi%kwchoi12@netvigator.com%1,                     YM%Sung%mandysung@cuhk.edu.hk%1,                     Paul KS%Chan%paulkschan@cuhk.edu.hk%1,                     YK%Tong%yktong@cuhk.edu.hk%1,                     ST%Lai%laist@ha.org.hk%1,                     WC%Yu%yuwc@ha.org.hk%1,                     Owen%Tsang%tyotsang@hotmail.com%1,                     YM Dennis%Lo%loym@cuhk.edu.hk%1]"
$ws.Cells.Item(6, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(7, 5).Value = "[Emi%Hamano%NULL%1,                     Minako%Hijikata%NULL%0,                     Satoru%Itoyama%NULL%2,                     Tran%Quy%NULL%0,                     Nguyen Chi%Phi%NULL%0,                     Hoang Thuy%Long%NULL%0,                     Le Dang%Ha%NULL%0,                     Vo Van%Ban%NULL%0,                     Ikumi%Matsushita%NULL%0,                     Hideki%Yanai%NULL%0,                     Fumiko%Kirikae%NULL%2,                     Teruo%Kirikae%NULL%0,                     Tadatoshi%Kuratsuji%NULL%0,                     Takehiko%Sasazuki%NULL%0,                     Naoto%Keicho%nkeicho-tky@umin.ac.jp%0]"
$ws.Cells.Item(7, 9).Value = "_PMC_elsevier_CROSSREF"

$ws.Cells.Item(8, 5).Value = "[W. K. Eddie%Ip%NULL%1,                     Kwok Hung%Chan%NULL%1,                     Helen K. W.%Law%NULL%1,                     Gloria H. W.%Tso%NULL%1,                     Eric K. P.%Kong%NULL%1,                     Wilfred H. S.%Wong%NULL%1,                     Yuk Fai%To%NULL%1,                     Raymond W. H.%Yung%NULL%1,                     Eudora Y.%Chow%NULL%1,                     Ka Leung%Au%NULL%2,                     Eric Y. T.%Chan%NULL%2,                     Wilina%Lim%NULL%3,                     Jens C.%Jensenius%NULL%1,                     Malcolm W.%Turner%NULL%1,                     J. S. Malik%Peiris%NULL%1,                     Yu Lung%Lau%NULL%3]"
$ws.Cells.Item(8, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(9, 5).Value = "[Satoru%Itoyama%NULL%0,                     Naoto%Keicho%nkeicho-tky@umin.ac.jp%0,                     Minako%Hijikata%NULL%0,                     Tran%Quy%NULL%0,                     Nguyen Chi%Phi%NULL%0,                     Hoang Thuy%Long%NULL%0,                     Le Dang%Ha%NULL%0,                     Vo Van%Ban%NULL%0,                     Ikumi%Matsushita%NULL%0,                     Hideki%Yanai%NULL%0,                     Fumiko%Kirikae%NULL%0,                     Teruo%Kirikae%NULL%0,                     Tadatoshi%Kuratsuji%NULL%0,                     Takehiko%Sasazuki%NULL%0]"
$ws.Cells.Item(9, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(10, 5).Value = "[F. F.%Yuan%NULL%1,                     J.%Tanner%NULL%1,                     P. K. S.%Chan%NULL%1,                     S.%Biffin%NULL%1,                     W. B.%Dyer%NULL%1,                     A. F.%Geczy%NULL%1,                     J. W.%Tang%NULL%1,                     D. S. C.%Hui%NULL%1,                     J. J. Y.%Sung%NULL%1,                     J. S.%Sullivan%NULL%1]"
$ws.Cells.Item(10, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(11, 5).Value = "[Hongxing%Zhang%NULL%1,                     Gangqiao%Zhou%NULL%2,                     Lianteng%Zhi%NULL%1,                     Hao%Yang%NULL%1,                     Yun%Zhai%NULL%2,                     Xiaojia%Dong%NULL%1,                     Xiumei%Zhang%NULL%1,                     Xue%Gao%NULL%1,                     Yunping%Zhu%NULL%1,                     Fuchu%He%NULL%2]"
$ws.Cells.Item(11, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(12, 5).Value = "[Vera S F%Chan%NULL%1,                     Kelvin Y K%Chan%NULL%1,                     Yongxiong%Chen%NULL%1,                     Leo L M%Poon%NULL%1,                     Annie N Y%Cheung%NULL%1,                     Bojian%Zheng%NULL%1,                     Kwok-Hung%Chan%NULL%1,                     William%Mak%NULL%1,                     Hextan Y S%Ngan%NULL%1,                     Xiaoning%Xu%NULL%1,                     Gavin%Screaton%NULL%1,                     Paul K H%Tam%NULL%1,                     Jonathan M%Austyn%NULL%1,                     Li-Chong%Chan%NULL%2,                     Shea-Ping%Yip%NULL%2,                     Malik%Peiris%NULL%2,                     Ui-Soon%Khoo%uskhoo@pathology.hku.hk%4,                     Chen-Lung S%Lin%clin@hkucc.hku.hk%1]"
$ws.Cells.Item(12, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(13, 5).Value = "[Wai Po% Chong%h9820905@graduate.hku.hk%1,                     WK Eddie%Ip%WIP@partners.org%1,                     Gloria Hoi Wan%Tso%h98083405@hkusua.hku.hk%1,                     Man Wai%Ng%ivy_natsu@yahoo.com.hk%1,                     Wilfred Hing Sang%Wong%whswong@hkucc.hku.hk%2,                     Helen Ka Wai%Law%hkwlaw@hkucc.hku.hk%2,                     Raymond WH%Yung%rwhyung@ha.org.hk%2,                     Eudora Y%Chow%chowe@ha.org.hk%2,                     KL%Au%klau@ha.org.hk%1,                     Eric YT%Chan%eytchan@ha.org.hk%2,                     Wilina%Lim%wllim@pacific.net.hk%0,                     JS Malik%Peiris%malik@hkucc.hku.hk%2,                     Yu Lung%Lau%lauylung@hkucc.hku.hk%0]"
$ws.Cells.Item(13, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(14, 5).Value = "[Jing%He%jinghe2004@gmail.com%1,                     Dan%Feng%fddd@263.net%1,                     Sake J%de Vlas%s.devlas@erasmusmc.nl%1,                     Hongwei%Wang%whw7809@yahoo.com.cn%1,                     Arnaud%Fontanet%fontanet@pasteur.fr%1,                     Panhe%Zhang%juntianz@yahoo.com.cn%1,                     Sabine%Plancoulaine%plancoulaine@necker.fr%1,                     Fang%Tang%tf4065@sina.com%1,                     Lin%Zhan%zhanzhan_30@yahoo.com.cn%1,                     Hong%Yang%anni_yang@yahoo.com%1,                     Tianbao%Wang%tianbaowang2005@yahoo.com.cn%1,                     Jan H%Richardus%j.richardus@erasmusmc.nl%1,                     J Dik F%Habbema%j.d.f.habbema@erasmusmc.nl%1,                     Wuchun%Cao%caowc@nic.bmi.ac.cn%1]"
$ws.Cells.Item(14, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(15, 5).Value = "[Wei-Ju%Chen%NULL%1,                     Jyh-Yuan%Yang%NULL%1,                     Jih-Hui%Lin%NULL%1,                     Cathy S. J.%Fann%NULL%1,                     Valeriy%Osyetrov%NULL%1,                     Chwan-Chuen%King%NULL%1,                     Yi-Ming Arthur%Chen%NULL%1,                     Hsiao-Ling%Chang%NULL%1,                     Hung-Wei%Kuo%NULL%1,                     Fong%Liao%NULL%1,                     Mei-Shang%Ho%homs@ibms.sinica.edu.tw%1]"
$ws.Cells.Item(15, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(16, 5).Value = "[Yi-Ming Arthur%Chen%xref no email%1,            Shu-Yuan%Liang%xref no email%1,            Yi-Ping%Shih%xref no email%1,            Chia-Yen%Chen%xref no email%1,            Yuan-Ming%Lee%xref no email%1,            Ling%Chang%xref no email%1,            Shiao-Ying%Jung%xref no email%1,            Mei-Shang%Ho%xref no email%1,            Kung-Yee%Liang%xref no email%1,            Hour-Young%Chen%xref no email%1,            Yu-Jiun%Chan%xref no email%1,            Da-Chen%Chu%xref no email%1]"
$ws.Cells.Item(16, 9).Value = "_CROSSREF"

$ws.Cells.Item(17, 5).Value = "[Kelvin Y. K.%Chan%NULL%1,                     Johannes C. Y.%Ching%NULL%1,                     M. S.%Xu%NULL%1,                     Annie N. Y.%Cheung%NULL%1,                     Shea-Ping%Yip%NULL%0,                     Loretta Y. C.%Yam%NULL%1,                     Sik-To%Lai%NULL%2,                     Chung-Ming%Chu%NULL%2,                     Andrew T. Y.%Wong%NULL%1,                     You-Qiang%Song%NULL%1,                     Fang-Ping%Huang%NULL%1,                     Wei%Liu%NULL%1,                     P. H.%Chung%NULL%1,                     G. M.%Leung%NULL%1,                     Eudora Y. D.%Chow%NULL%1,                     Eric Y. T.%Chan%NULL%0,                     Jane C. K.%Chan%NULL%1,                     Hextan%Ngan%NULL%1,                     Paul%Tam%NULL%1,                     Li-Chong%Chan%NULL%0,                     Pak%Sham%NULL%1,                     Vera S. F.%Chan%NULL%1,                     Malik%Peiris%NULL%0,                     Steve C. L.%Lin%NULL%1,                     Ui-Soon%Khoo%uskhoo@pathology.hku.hk%0]"
$ws.Cells.Item(17, 9).Value = "_PMC_CROSSREF"

$ws.Cells.Item(18, 5).Value = "[Man Wai%Ng%ivy_natsu@yahoo.com.hk%1,                     Gangqiao%Zhou%zhougq@chgb.org.cn%0,                     Wai Po%Chong%h9820905@graduate.hku.hk%1,                     Loretta Wing Yan%Lee%h0226204@hkusua.hku.hk%1,                     Helen Ka Wai%Law%hkwlaw@hkucc.hku.hk%0,                     Hongxing%Zhang%zhanghx08@126.com%1,                     Wilfred Hing Sang%Wong%whswong@hkucc.hku.hk%0,                     Susanna Fung Shan%Fok%sfsfok@hkusua.hku.hk%1,                     Yun%Zhai%zhaiyun077@tom.com%0,                     Raymond WH%Yung%rwhyung@ha.org.hk%0,                     Eudora Y%Chow%chowe@ha.org.hk%0,                     Ka Leung%Au%klau@ha.org.hk%0,                     Eric YT%Chan%eytchan@ha.org.hk%0,                     Wilina%Lim%wllim@pacific.net.hk%0,                     JS Malik%Peiris%malik@hkucc.hku.hk%0,                     Fuchu%He%hefc@nic.bmi.ac.cn%0,                     Yu Lung%Lau%lauylung@hkucc.hku.hk%0]"
$ws.Cells.Item(18, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(19, 5).Value = "[Fang F.%Yuan%xref no email%1,            Ingrid%Boehm%xref no email%1,            Paul K. S.%Chan%xref no email%1,            Katherine%Marks%xref no email%1,            Julian W.%Tang%xref no email%1,            David S. C.%Hui%xref no email%1,            Joseph J. Y.%Sung%xref no email%1,            Wayne B.%Dyer%xref no email%1,            Andrew F.%Geczy%xref no email%1,            John S.%Sullivan%xref no email%1]"
$ws.Cells.Item(19, 9).Value = "_CROSSREF"

$ws.Cells.Item(20, 5).Value = "[H%Li%xref no email%1,            N L-S%Tang%xref no email%1,            P K-S%Chan%xref no email%1,            C-Y%Wang%xref no email%1,            D S-C%Hui%xref no email%1,            C%Luk%xref no email%1,            R%Kwok%xref no email%1,            W%Huang%xref no email%1,            J J-Y%Sung%xref no email%1,            Q-P%Kong%xref no email%1,            Y-P%Zhang%xref no email%1]"
$ws.Cells.Item(20, 9).Value = "_CROSSREF"

$ws.Cells.Item(21, 5).Value = "[Shixin%Wang%Wshx-001@163.com%1,                     Maoti%Wei%weimaoti@yahoo.com.cn%1,                     Yi%Han%halanhan2001sh@163.com%1,                     Keju%Zhang%sinong77@sina.com%1,                     Li%He%wujinhl@163.com%1,                     Zhen%Yang%yzdtchina@yahoo.com.cn%1,                     Bing%Su%tjsubin553@sina.com%1,                     Zhilun%Zhang%zzl718@sina.com.cn%1,                     Yilan%Hu%yougu@yahoo.com.cn%1,                     Wuli%Hui%huiwuli@163.com%1]"
$ws.Cells.Item(21, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(22, 5).Value = "[Naoto%Keicho%nkeicho-tky@umin.ac.jp%1,                     Satoru%Itoyama%NULL%0,                     Koichi%Kashiwase%NULL%1,                     Nguyen Chi%Phi%NULL%0,                     Hoang Thuy%Long%NULL%0,                     Le Dang%Ha%NULL%0,                     Vo Van%Ban%NULL%0,                     Bach Khanh%Hoa%NULL%1,                     Nguyen Thi Le%Hang%NULL%1,                     Minako%Hijikata%NULL%0,                     Shinsaku%Sakurada%NULL%1,                     Masahiro%Satake%NULL%1,                     Katsushi%Tokunaga%NULL%1,                     Takehiko%Sasazuki%NULL%0,                     Tran%Quy%NULL%0]"
$ws.Cells.Item(22, 9).Value = "_PMC_elsevier_CROSSREF"

$ws.Cells.Item(23, 5).Value = "[Yan%Wang%001wangyan@sina.com%1,                     Jiangwei%Yan%jiangweiyanbj@yahoo.com.cn%1,                     Yuling%Shi%sane222@163.com%1,                     Ping%Li%liping@nic.bmi.ac.cn%1,                     Chuanxuan%Liu%renew7306_cn@sina.com%1,                     Qingjun%Ma%maqingjun@sina.com%1,                     Ruifu%Yang%zhangf@big.ac.cn%1,                     Xiaoyi%Wang%happybettylx@hotmail.com%1,                     Lina%zhu%luckina@sina.com%1,                     Xiao%Yang%gaowei@sina.com%1,                     Cheng%Cao%caoc@nic.bmi.ac.cn%1]"
$ws.Cells.Item(23, 9).Value = "_PMC_Springer_CROSSREF"

$ws.Cells.Item(24, 5).Value = "[Kelvin Yuen Kwong%Chan%NULL%1,                     Mei-Shu%Xu%NULL%2,                     Johannes Chi Yun%Ching%NULL%1,                     Thomas Man Kit%So%NULL%1,                     Sik-To%Lai%NULL%0,                     Chung-Ming%Chu%NULL%0,                     Loretta Y.C.%Yam%NULL%1,                     Andrew T.Y.%Wong%NULL%1,                     Pui Hong%Chung%NULL%1,                     Vera Sau Fong%Chan%NULL%1,                     Chen Lung Steve%Lin%NULL%1,                     Pak Chung%Sham%NULL%1,                     Gabriel M.%Leung%NULL%2,                     Joseph S.M.%Peiris%NULL%1,                     Ui-Soon%Khoo%uskhoo@pathology.hku.hk%0]"
$ws.Cells.Item(24, 9).Value = "_PMC_elsevier_CROSSREF"

$ws.Cells.Item(25, 5).Value = "[Johannes Chi-Yun%Ching%NULL%1,                     Kelvin Yuen Kwong%Chan%NULL%1,                     Eric Hing Leung%Lee%NULL%1,                     Mei-Shu%Xu%NULL%0,                     Campbell Kam Po%Ting%NULL%1,                     Thomas M. K.%So%NULL%1,                     Pak C.%Sham%NULL%1,                     Gabriel M.%Leung%NULL%0,                     Joseph S. M.%Peiris%NULL%1,                     Ui-Soon%Khoo%uskhoo@pathology.hku.hk%0]"
$ws.Cells.Item(25, 9).Value = "_PMC_CROSSREF"
